# Swap the presentation's theme palette from the "Integral" (Red Violet)
# scheme to the stock "Office Theme" (Office) scheme.
#
# PowerPoint's ColorScheme.Colors(n).RGB setter expects a packed BGR
# integer (the classic VBA "RGB()" long: 0x00BBGGRR), so convert each
# target hex "RRGGBB" triad accordingly.
function HexToComRgb($hex) {
    $v = [Convert]::ToInt32($hex, 16)
    $r = ($v -shr 16) -band 0xFF
    $g = ($v -shr 8) -band 0xFF
    $b = $v -band 0xFF
    return $b * 65536 + $g * 256 + $r
}

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.ColorScheme

# Office Theme color scheme, in clrScheme document order:
# dk1, lt1, dk2, lt2, accent1..accent6, hlink, folHlink
$officeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

for ($i = 0; $i -lt $officeColors.Count; $i++) {
    $cs.Colors($i + 1).RGB = HexToComRgb($officeColors[$i])
}
